$wb = $excel.ActiveWorkbook

# --- Layer0 worksheet data (Bias / Weight0..Weight4 rows) ---
$ws1 = $wb.Worksheets.Item("Layer0")

$layer0Data = @(
    @(-4.651968112421115,0.3199272236125825,-1.213589177504441,0.05882924384248901,-2.721855465877633,3.687301547932157,-3.154517396286086,-8.973328325686738,0.7887606495005158,-0.8773487503236762,0.7139641943800322,5.100921333029097,-3.221416203891608,4.455733949134388,6.104686149709434,-0.5110331822840476),
    @(-1.355289232618496,1.498178373125511,3.076712849820554,0.5451097286303356,-1.671326240878229,-2.285656069509277,-0.3740969033482252,-2.234930095312851,5.170240301407581,0.8898644087355447,0.3096916883932206,-3.625325201393851,-3.270514088977783,0.779487555070089,1.01566080836038,-2.941699634952081),
    @(3.78613942150798,2.191256352398883,0.7689871978256049,8.253650831768791,-4.085025275706212,-1.962804743283969,-2.718137087113132,2.212742426810743,1.405612970279424,-3.160244109890405,1.76448642389058,-3.884841821281861,-2.763162736267853,-8.036676673404795,-1.079354068257501,-0.6454604972423791),
    @(-1.595020618210012,-0.6095588971095882,-1.969264812872284,-2.219013037302394,1.650480393646374,2.832849948792493,1.981834075831763,-0.3632443234000373,-0.8523474583975876,3.07124312264388,0.9602951137130332,-1.733151876420568,1.780653230389478,0.8699209591603677,-6.766446181995086,2.075145973636541),
    @(-0.778565950607754,6.246925885220882,-4.040357870808652,0.1949037966116655,2.474669381968873,4.960133848788453,-2.512214490642664,-8.930551866512268,-1.063690816624645,-4.293550632813429,-0.6653553823763896,1.594813428431989,4.061310393999485,1.20736868850456,0.1028660952319049,-0.246291740545291),
    @(2.348703619816927,-0.9869239512761522,1.490154614888753,0.819513033332126,0.3407076974492702,-3.200643351501788,-2.468421587716509,0.09388305768003274,-6.404304572117677,-3.623237482413045,-3.354463130684142,6.694239313449123,2.676477416826273,1.931556765593156,2.978890138115263,-6.558043444040501)
)

for ($r = 0; $r -lt $layer0Data.Length; $r++) {
    $rowVals = $layer0Data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws1.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}

# --- Layer1 worksheet data (Bias / Weight0..Weight15 rows) ---
$ws2 = $wb.Worksheets.Item("Layer1")

$layer1Data = @(
    @(-3.67877393229235,-4.559229601873704,-1.795812826129106,-4.493073084950532,-1.816118240887888),
    @(1.37875237629632,-2.383565229616215,1.564534506942169,-2.27221244343842,-5.486375579984861),
    @(3.13492262151221,-2.278479595797708,-4.975181069484655,1.807849160074941,-0.180513379099211),
    @(-4.029155944660562,-4.200831107853952,3.540578668196304,1.172162859698816,-3.613748432004836),
    @(-2.616878907824075,-5.352822294509116,3.963673270340919,-6.694497631964569,4.471729198419041),
    @(4.391405553220466,3.379695874319064,-0.5529734841351881,-5.599077270787865,-5.708626251008716),
    @(0.2603714250372132,-1.206259173061807,-4.867326906775545,-5.401514226454959,4.148691235148113),
    @(-3.631416560820565,3.522903267833373,0.02359518056822767,2.147193914060618,-3.899377559726754),
    @(-1.579003904740546,-3.600647863225803,4.756874477663022,-6.665633419237031,-1.431661294217504),
    @(-5.733641121959042,0.2187783419972568,-3.616510945878727,4.592579307285932,-3.761998626554248),
    @(-2.78541731107352,3.311282012717154,-1.445455048731638,2.290788933144136,-2.422485094344136),
    @(-4.828790099621407,1.596151895787254,-0.03277843478382848,-0.9962989724158545,1.521811805947623),
    @(2.861809728919606,-3.957042194711581,-4.006528987355455,5.373336847488569,-2.289803944693718),
    @(2.693541888521888,4.171058537088271,-2.08007914907212,-7.618266479437511,-2.010353212917932),
    @(0.9844179201164096,2.23450148900825,-4.331071536478865,2.147253319961016,-2.406868734515847),
    @(0.4618054340097657,-8.117856574191403,-0.2747218844424215,5.358480436770769,0.1567277929826223),
    @(-6.315692667296799,0.5682751228517877,0.3031953747597756,-3.616054561302839,6.828743496034728)
)

for ($r = 0; $r -lt $layer1Data.Length; $r++) {
    $rowVals = $layer1Data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws2.Cells.Item($r + 2, $c + 2).Value = $rowVals[$c]
    }
}
